$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EV Away win")

# Insert a new row at position 12, shifting the existing row 12 (and below) down to row 13
$ws.Rows.Item(12).Insert()

$ws.Cells.Item(12, 1).Value = "20-05-2025 19:00"
$ws.Cells.Item(12, 2).Value = "NORWAY"
$ws.Cells.Item(12, 3).Value = "NM CUPEN"
$ws.Cells.Item(12, 4).Value = "Bryne - KFUM Oslo"
$ws.Cells.Item(12, 5).Value = 51.7
$ws.Cells.Item(12, 6).Value = 2.8
$ws.Cells.Item(12, 7).Value = 0.45
